$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove AA183 and AD183 (now-empty marker cells dropped from row 183) ---
$ws.Range("AA183").ClearContents()
$ws.Range("AD183").ClearContents()

# --- Row 184 ---
$ws.Range("A184:Z184").NumberFormat = "@"
$ws.Range("AB184:AC184").NumberFormat = "@"
$ws.Range("AE184:BI184").NumberFormat = "@"

$ws.Range("A184").Value = "2022-07-26 14:24:35"
$ws.Range("B184").Value = "26.6"
$ws.Range("C184").Value = "27.2"
$ws.Range("D184").Value = "25.8"
$ws.Range("E184").Value = "."
$ws.Range("F184").Value = "25.6"
$ws.Range("G184").Value = "."
$ws.Range("H184").Value = "53"
$ws.Range("I184").Value = "."
$ws.Range("J184").Value = "20"
$ws.Range("K184").Value = "20"
$ws.Range("L184").Value = "0"
$ws.Range("M184").Value = "0"
$ws.Range("N184").Value = "20"
$ws.Range("O184").Value = "20"
$ws.Range("P184").Value = "."
$ws.Range("Q184").Value = "."
$ws.Range("R184").Value = "0"
$ws.Range("S184").Value = "0"
$ws.Range("T184").Value = "0"
$ws.Range("U184").Value = "0"
$ws.Range("V184").Value = "100"
$ws.Range("W184").Value = "72"
$ws.Range("X184").Value = "."
$ws.Range("Y184").Value = "."
$ws.Range("Z184").Value = "12.61"
$ws.Range("AB184").Value = "7002200"
$ws.Range("AC184").Value = "2200"
$ws.Range("AE184").Value = "26.7"
$ws.Range("AF184").Value = "27.3"
$ws.Range("AG184").Value = "25.7"
$ws.Range("AH184").Value = "20"
$ws.Range("AI184").Value = "20"
$ws.Range("AJ184").Value = "20"
$ws.Range("AK184").Value = "50"
$ws.Range("AL184").Value = "50"
$ws.Range("AM184").Value = "72"
$ws.Range("AN184").Value = "."
$ws.Range("AO184").Value = "0"
$ws.Range("AP184").Value = "10"
$ws.Range("AQ184").Value = "0"
$ws.Range("AR184").Value = "0.15"
$ws.Range("AS184").Value = "0.16"
$ws.Range("AT184").Value = "0.21"
$ws.Range("AU184").Value = "11.23"
$ws.Range("AV184").Value = "3584.69"
$ws.Range("AW184").Value = "0.00"
$ws.Range("AX184").Value = "0.00"
$ws.Range("AY184").Value = "2913.20"
$ws.Range("AZ184").Value = "0.00"
$ws.Range("BA184").Value = "4.68"
$ws.Range("BB184").Value = "10331.15"
$ws.Range("BC184").Value = "19.0"
$ws.Range("BD184").Value = "."
$ws.Range("BE184").Value = "."
$ws.Range("BF184").Value = "53"
$ws.Range("BG184").Value = "."
$ws.Range("BH184").Value = "53"
$ws.Range("BI184").Value = "149176580"

# --- Row 185 ---
$ws.Range("A185:BI185").NumberFormat = "@"

$ws.Range("A185").Value = "2022-07-26 14:24:49"
$ws.Range("B185").Value = "26.6"
$ws.Range("C185").Value = "27.1"
$ws.Range("D185").Value = "25.7"
$ws.Range("E185").Value = "."
$ws.Range("F185").Value = "25.6"
$ws.Range("G185").Value = "."
$ws.Range("H185").Value = "53"
$ws.Range("I185").Value = "."
$ws.Range("J185").Value = "20"
$ws.Range("K185").Value = "20"
$ws.Range("L185").Value = "0"
$ws.Range("M185").Value = "0"
$ws.Range("N185").Value = "20"
$ws.Range("O185").Value = "20"
$ws.Range("P185").Value = "."
$ws.Range("Q185").Value = "."
$ws.Range("R185").Value = "0"
$ws.Range("S185").Value = "0"
$ws.Range("T185").Value = "0"
$ws.Range("U185").Value = "0"
$ws.Range("V185").Value = "100"
$ws.Range("W185").Value = "72"
$ws.Range("X185").Value = "."
$ws.Range("Y185").Value = "."
$ws.Range("Z185").Value = "12.61"
$ws.Range("AA185").Value = "'"
$ws.Range("AB185").Value = "7002200"
$ws.Range("AC185").Value = "2200"
$ws.Range("AD185").Value = "'"
$ws.Range("AE185").Value = "26.5"
$ws.Range("AF185").Value = "27.2"
$ws.Range("AG185").Value = "25.7"
$ws.Range("AH185").Value = "20"
$ws.Range("AI185").Value = "20"
$ws.Range("AJ185").Value = "20"
$ws.Range("AK185").Value = "50"
$ws.Range("AL185").Value = "50"
$ws.Range("AM185").Value = "72"
$ws.Range("AN185").Value = "."
$ws.Range("AO185").Value = "0"
$ws.Range("AP185").Value = "10"
$ws.Range("AQ185").Value = "0"
$ws.Range("AR185").Value = "0.15"
$ws.Range("AS185").Value = "0.16"
$ws.Range("AT185").Value = "0.21"
$ws.Range("AU185").Value = "11.23"
$ws.Range("AV185").Value = "3584.69"
$ws.Range("AW185").Value = "0.00"
$ws.Range("AX185").Value = "0.00"
$ws.Range("AY185").Value = "2913.20"
$ws.Range("AZ185").Value = "0.00"
$ws.Range("BA185").Value = "4.68"
$ws.Range("BB185").Value = "10331.15"
$ws.Range("BC185").Value = "19.0"
$ws.Range("BD185").Value = "."
$ws.Range("BE185").Value = "."
$ws.Range("BF185").Value = "53"
$ws.Range("BG185").Value = "."
$ws.Range("BH185").Value = "53"
$ws.Range("BI185").Value = "149176580"
